# Strip text for tabular export
# Trailing newlines/whitespace in the "text" column of the exported
# requirements table are stripped so each cell holds clean, trimmed text
# (matches the shared-string cleanup: "Hello, world!\n" -> "Hello, world!",
# etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# Locate the "text" column by scanning the header row (row 1).
$headerRow = 1
$textCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item($headerRow, $c).Value2
    if ($header -eq "text") {
        $textCol = $c
    }
}

if ($textCol -gt 0) {
    for ($r = $headerRow + 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, $textCol)
        $value = $cell.Value2
        if ($value -ne $null) {
            $cell.Value2 = $value.ToString().Trim()
        }
    }
}
